# Populate the "HALO EXCHANGE TIMES" timing results (Time 1-5, columns X:AB)
# for the two data blocks (rows 93-111 and 115-132). Column AC already holds
# a shared "=SUM(X#:AB#)/5" formula per row, so it recalculates automatically
# once the Time 1-5 cells are populated - no need to touch it directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("X93").Value = 0.0043
$ws.Range("Y93").Value = 0.0049760000000000004
$ws.Range("Z93").Value = 0.0043
$ws.Range("AA93").Value = 0.0043509999999999998
$ws.Range("AB93").Value = 0.0043140000000000001

$ws.Range("X94").Value = 0.0054679999999999998
$ws.Range("Y94").Value = 0.0067010000000000004
$ws.Range("Z94").Value = 0.0064650000000000003
$ws.Range("AA94").Value = 0.006084
$ws.Range("AB94").Value = 0.0054510000000000001

$ws.Range("X95").Value = 0.006881
$ws.Range("Y95").Value = 0.0060939999999999996
$ws.Range("Z95").Value = 0.0059069999999999999
$ws.Range("AA95").Value = 0.0069369999999999996
$ws.Range("AB95").Value = 0.0063049999999999998

$ws.Range("X96").Value = 0.0056990000000000001
$ws.Range("Y96").Value = 0.0058040000000000001
$ws.Range("Z96").Value = 0.0056730000000000001
$ws.Range("AA96").Value = 0.0059100000000000003
$ws.Range("AB96").Value = 0.004463

$ws.Range("X97").Value = 0.0063509999999999999
$ws.Range("Y97").Value = 0.0068950000000000001
$ws.Range("Z97").Value = 0.0069199999999999999
$ws.Range("AA97").Value = 0.0065189999999999996
$ws.Range("AB97").Value = 0.0059090000000000002

$ws.Range("X98").Value = 0.0089990000000000001
$ws.Range("Y98").Value = 0.004607
$ws.Range("Z98").Value = 0.005522
$ws.Range("AA98").Value = 0.0072610000000000001
$ws.Range("AB98").Value = 0.0066119999999999998

$ws.Range("X99").Value = 0.0044749999999999998
$ws.Range("Y99").Value = 0.0060280000000000004
$ws.Range("Z99").Value = 0.006215
$ws.Range("AA99").Value = 0.0065750000000000001
$ws.Range("AB99").Value = 0.0045269999999999998

$ws.Range("X100").Value = 0.10738
$ws.Range("Y100").Value = 0.0075799999999999999
$ws.Range("Z100").Value = 0.0046090000000000002
$ws.Range("AA100").Value = 0.095135999999999998
$ws.Range("AB100").Value = 0.068278000000000005

$ws.Range("X101").Value = 0.0087030000000000007
$ws.Range("Y101").Value = 0.0046179999999999997
$ws.Range("Z101").Value = 0.0044970000000000001
$ws.Range("AA101").Value = 0.005032
$ws.Range("AB101").Value = 0.0049880000000000002

$ws.Range("X102").Value = 0.0068019999999999999
$ws.Range("Y102").Value = 0.172598
$ws.Range("Z102").Value = 0.110495
$ws.Range("AA102").Value = 0.162693
$ws.Range("AB102").Value = 0.165802

$ws.Range("X103").Value = 0.27704000000000001
$ws.Range("Y103").Value = 0.12509799999999999
$ws.Range("Z103").Value = 0.0099830000000000006
$ws.Range("AA103").Value = 0.0077640000000000001
$ws.Range("AB103").Value = 0.0066319999999999999

$ws.Range("X104").Value = 0.27704000000000001
$ws.Range("Y104").Value = 0.12509799999999999
$ws.Range("Z104").Value = 0.0099830000000000006
$ws.Range("AA104").Value = 0.0077640000000000001
$ws.Range("AB104").Value = 0.0066319999999999999

$ws.Range("X105").Value = 0.177564
$ws.Range("Y105").Value = 0.24526800000000001
$ws.Range("Z105").Value = 0.206123
$ws.Range("AA105").Value = 0.29887200000000003
$ws.Range("AB105").Value = 0.41246100000000002

$ws.Range("X106").Value = 0.21778400000000001
$ws.Range("Y106").Value = 0.024847999999999999
$ws.Range("Z106").Value = 0.010624
$ws.Range("AA106").Value = 0.22184699999999999
$ws.Range("AB106").Value = 0.36779099999999998

$ws.Range("X107").Value = 0.050848999999999998
$ws.Range("Y107").Value = 0.26322499999999999
$ws.Range("Z107").Value = 0.18569099999999999
$ws.Range("AA107").Value = 0.23388800000000001
$ws.Range("AB107").Value = 0.20404800000000001

$ws.Range("X108").Value = 0.113939
$ws.Range("Y108").Value = 0.183503
$ws.Range("Z108").Value = 0.12995599999999999
$ws.Range("AA108").Value = 0.28311599999999998
$ws.Range("AB108").Value = 0.14461199999999999

$ws.Range("X109").Value = 0.47208
$ws.Range("Y109").Value = 0.13689299999999999
$ws.Range("Z109").Value = 0.069712999999999997
$ws.Range("AA109").Value = 0.11792999999999999
$ws.Range("AB109").Value = 0.091983999999999996

$ws.Range("X110").Value = 0.56765100000000002
$ws.Range("Y110").Value = 0.013899
$ws.Range("Z110").Value = 0.012274
$ws.Range("AA110").Value = 0.029822999999999999
$ws.Range("AB110").Value = 0.014338999999999999

$ws.Range("X111").Value = 0.62446299999999999
$ws.Range("Y111").Value = 0.047230000000000001
$ws.Range("Z111").Value = 0.047632000000000001
$ws.Range("AA111").Value = 0.10287300000000001
$ws.Range("AB111").Value = 0.023654000000000001

$ws.Range("X115").Value = 0.0095589999999999998
$ws.Range("Y115").Value = 0.011979999999999999
$ws.Range("Z115").Value = 0.013820000000000001
$ws.Range("AA115").Value = 0.017468999999999998
$ws.Range("AB115").Value = 0.016319

$ws.Range("X116").Value = 0.021232999999999998
$ws.Range("Y116").Value = 0.021315000000000001
$ws.Range("Z116").Value = 0.010711999999999999
$ws.Range("AA116").Value = 0.015858000000000001
$ws.Range("AB116").Value = 0.024036999999999999

$ws.Range("X117").Value = 0.02649
$ws.Range("Y117").Value = 0.013974
$ws.Range("Z117").Value = 0.012239
$ws.Range("AA117").Value = 0.012086
$ws.Range("AB117").Value = 0.012194999999999999

$ws.Range("X118").Value = 0.015247
$ws.Range("Y118").Value = 0.027316
$ws.Range("Z118").Value = 0.023779000000000002
$ws.Range("AA118").Value = 0.015633000000000001
$ws.Range("AB118").Value = 0.015491

$ws.Range("X119").Value = 0.025208999999999999
$ws.Range("Y119").Value = 0.039017999999999997
$ws.Range("Z119").Value = 0.014037000000000001
$ws.Range("AA119").Value = 0.030346999999999999
$ws.Range("AB119").Value = 0.017994

$ws.Range("X120").Value = 0.047135999999999997
$ws.Range("Y120").Value = 0.028195999999999999
$ws.Range("Z120").Value = 0.025918
$ws.Range("AA120").Value = 0.029340999999999999
$ws.Range("AB120").Value = 0.029142000000000001

$ws.Range("X121").Value = 0.080348000000000003
$ws.Range("Y121").Value = 0.028313999999999999
$ws.Range("Z121").Value = 0.033370999999999998
$ws.Range("AA121").Value = 0.032885999999999999
$ws.Range("AB121").Value = 0.031556000000000001

$ws.Range("X122").Value = 0.018183000000000001
$ws.Range("Y122").Value = 0.031674000000000001
$ws.Range("Z122").Value = 0.025736999999999999
$ws.Range("AA122").Value = 0.031261999999999998
$ws.Range("AB122").Value = 0.031168999999999999

$ws.Range("X123").Value = 0.012538000000000001
$ws.Range("Y123").Value = 0.013513000000000001
$ws.Range("Z123").Value = 0.013903
$ws.Range("AA123").Value = 0.020648
$ws.Range("AB123").Value = 0.017177999999999999

$ws.Range("X124").Value = 0.020282000000000001
$ws.Range("Y124").Value = 0.013311
$ws.Range("Z124").Value = 0.025440000000000001
$ws.Range("AA124").Value = 0.028660999999999999
$ws.Range("AB124").Value = 0.020795000000000001

$ws.Range("X125").Value = 0.033696999999999998
$ws.Range("Y125").Value = 0.044817999999999997
$ws.Range("Z125").Value = 0.013408
$ws.Range("AA125").Value = 0.046197000000000002
$ws.Range("AB125").Value = 0.018249999999999999

$ws.Range("X126").Value = 0.079344999999999999
$ws.Range("Y126").Value = 0.013127
$ws.Range("Z126").Value = 0.05151
$ws.Range("AA126").Value = 0.059158000000000002
$ws.Range("AB126").Value = 0.014533000000000001

$ws.Range("X127").Value = 0.018109
$ws.Range("Y127").Value = 0.038924
$ws.Range("Z127").Value = 0.037435000000000003
$ws.Range("AA127").Value = 0.011638000000000001
$ws.Range("AB127").Value = 0.050946999999999999

$ws.Range("X128").Value = 0.013358
$ws.Range("Y128").Value = 0.013207999999999999
$ws.Range("Z128").Value = 0.078900999999999999
$ws.Range("AA128").Value = 0.089417999999999997
$ws.Range("AB128").Value = 0.12805

$ws.Range("X129").Value = 0.027536000000000001
$ws.Range("Y129").Value = 0.012418
$ws.Range("Z129").Value = 0.10438
$ws.Range("AA129").Value = 0.030084
$ws.Range("AB129").Value = 0.020989000000000001

$ws.Range("X130").Value = 0.19158600000000001
$ws.Range("Y130").Value = 0.15037800000000001
$ws.Range("Z130").Value = 0.150231
$ws.Range("AA130").Value = 0.12665799999999999
$ws.Range("AB130").Value = 0.13945099999999999

$ws.Range("X131").Value = 0.10864500000000001
$ws.Range("Y131").Value = 0.170991
$ws.Range("Z131").Value = 0.131216
$ws.Range("AA131").Value = 0.17532500000000001
$ws.Range("AB131").Value = 0.186806

$ws.Range("X132").Value = 0.18029800000000001
$ws.Range("Y132").Value = 0.24496100000000001
$ws.Range("Z132").Value = 0.15149899999999999
$ws.Range("AA132").Value = 0.147511
$ws.Range("AB132").Value = 0.238089

# Update the sheet view's active cell/selection to match the saved state
# (Excel records the last selection in the worksheet XML on save).
[void]$ws.Range("AB133").Select()

